$d = $word.ActiveDocument

# The transcription contains editorial <corr>...</corr> markers (rendered
# as literal red Courier-New "<corr>"/"</corr>" text around a corrected
# letter, surrounded by normal black Arial prose text).
#
# 1) "... grosse<corr>t</corr>, mouille bien" -> "... grosset, mouille bien"
#    The <corr> / </corr> marker text is removed but the corrected letter
#    "t" itself is kept, merging back into the surrounding Arial text.
#
#    Find/Replace always takes on the formatting of the first character of
#    the match, so each marker is deleted starting from an adjoining plain
#    (Arial) character rather than from the marker itself; that way the
#    result keeps the correct Arial run formatting instead of inheriting
#    the red Courier New marker formatting.
$d.Content.Find.Execute("t</corr>", $true, $false, $false, $false, $false, $true, 1, $false, "t", 2) | Out-Null
$d.Content.Find.Execute("e<corr>", $true, $false, $false, $false, $false, $true, 1, $false, "e", 2) | Out-Null

# 2) "car il<corr>s</corr> se garderont humides bien un " ->
#    "car il se garderont humides bien un "
#    Here the whole <corr>s</corr> block (markers AND corrected letter "s")
#    is dropped. The match starts at "car il" (Arial), so the replacement
#    naturally keeps the correct Arial formatting.
$d.Content.Find.Execute("car il<corr>s</corr> se garderont humides bien un ", $true, $false, $false, $false, $false, $true, 1, $false, "car il se garderont humides bien un ", 2) | Out-Null
